$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70 -- this shifts the existing rows 70..91
# down to 71..92, matching the diff (which is a row-insert "weekly" update
# with the new observation placed right after the most-recent existing row).
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new weekly record.
$ws.Range("A70").Value = 9
$ws.Range("B70").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C70").Value = "Metropolitana"
$ws.Range("D70").Value = 45093
$ws.Range("E70").Value = 13
$ws.Range("F70").Value = 100112035
$ws.Range("G70").Value = "Bruselas (repollito)"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 52
$ws.Range("K70").Value = 17000
$ws.Range("L70").Value = 19000
$ws.Range("M70").Value = 18000
$ws.Range("N70").Value = "`$/malla 15 kilos"
$ws.Range("O70").Value = "Provincia de Quillota"
$ws.Range("P70").Value = 1200
$ws.Range("Q70").Value = 15
$ws.Range("R70").Value = "Hortaliza"
